# Actualización automática 2025-07-31 17:20:08
$wb = $excel.ActiveWorkbook

# --- Sheet 1: VENTAS POR GRUPO ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("E5").Value = 547.88
$ws1.Range("M5").Value = 8516.73

# --- Sheet 2: VENTA MENSUAL ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F5").Value = 9064.610000000001
$ws2.Range("F22").Value = 59179.46

# --- Sheet 3: CUMPLIMIENTO MENSUAL ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D4").Value = 547.88
$ws3.Range("E4").Value = -34.04895334066396
$ws3.Range("F4").Value = 1.066264881349682

$ws3.Range("D16").Value = 47937.16
$ws3.Range("E16").Value = -3670.920000000006
$ws3.Range("F16").Value = 1.082928208946592

$ws3.Range("D19").Value = 59179.46
$ws3.Range("E19").Value = 6198.537622917679
$ws3.Range("F19").Value = 0.9051892402904546
